$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6737.5
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 6737.5
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H138").Value = 4163.143
$ws.Range("J138").Value = 4351.449
$ws.Range("L138").Value = 13054.347
$ws.Range("N138").Value = -23334.347

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6034059
$ws.Range("I32").Value = 6415665.5
$ws.Range("K32").Value = 6415665.5
$ws.Range("M32").Value = -6415378.5
$ws.Range("H45").Value = 2488.077
$ws.Range("I45").Value = 2445
$ws.Range("K45").Value = 2445
$ws.Range("M45").Value = -2068
$ws.Range("H88").Value = 3122.111
$ws.Range("J88").Value = 2360
$ws.Range("L88").Value = 2360
$ws.Range("N88").Value = -3172
$ws.Range("H91").Value = 3122.111
$ws.Range("J91").Value = 2360
$ws.Range("L91").Value = 2360
$ws.Range("N91").Value = -5168
$ws.Range("H97").Value = 1438.1111
$ws.Range("I97").Value = 1146.7142
$ws.Range("J97").Value = 2458
$ws.Range("K97").Value = 1146.7142
$ws.Range("L97").Value = 2458
$ws.Range("M97").Value = -650.7141999999999
$ws.Range("N97").Value = -3450
$ws.Range("H122").Value = 3672.5
$ws.Range("I122").Value = 2103.7
$ws.Range("K122").Value = 6311.099999999999
$ws.Range("M122").Value = -3861.099999999999
$ws.Range("H132").Value = 5212.972
$ws.Range("I132").Value = 3833.6333
$ws.Range("K132").Value = 11500.8999
$ws.Range("M132").Value = -8970.8999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2894
$ws.Range("I86").Value = 2993.647
$ws.Range("K86").Value = 2993.647
$ws.Range("M86").Value = -1870.647
$ws.Range("H89").Value = 2894
$ws.Range("I89").Value = 2993.647
$ws.Range("K89").Value = 14968.235
$ws.Range("M89").Value = -9352.235000000001
$ws.Range("H99").Value = 11699.4
$ws.Range("I99").Value = 12454.889
$ws.Range("K99").Value = 12454.889
$ws.Range("M99").Value = -10956.889
$ws.Range("H107").Value = 1473.96
$ws.Range("I107").Value = 1401.5264
$ws.Range("K107").Value = 1401.5264
$ws.Range("M107").Value = 518.4736
$ws.Range("H130").Value = 97000
$ws.Range("J130").Value = 97000
$ws.Range("L130").Value = 97000
$ws.Range("N130").Value = -107040
$ws.Range("H134").Value = 29522.135
$ws.Range("I134").Value = 2564.0278
$ws.Range("K134").Value = 7692.0834
$ws.Range("M134").Value = -5157.0834

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 20471.143
$ws.Range("J43").Value = 20471.143
$ws.Range("L43").Value = 20471.143
$ws.Range("N43").Value = -20839.143
$ws.Range("H101").Value = 20471.143
$ws.Range("J101").Value = 20471.143
$ws.Range("L101").Value = 20471.143
$ws.Range("N101").Value = -26961.143
$ws.Range("H124").Value = 74661.664
$ws.Range("J124").Value = 74661.664
$ws.Range("L124").Value = 74661.664
$ws.Range("N124").Value = -79571.664
$ws.Range("H132").Value = 3580.1765
$ws.Range("I132").Value = 3616.4375
$ws.Range("K132").Value = 10849.3125
$ws.Range("M132").Value = -8319.3125

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 95.73684
$ws.Range("J2").Value = 276
$ws.Range("L2").Value = 1656
$ws.Range("N2").Value = -1882
$ws.Range("H134").Value = 16669967
$ws.Range("I134").Value = 20835458
$ws.Range("K134").Value = 62506374
$ws.Range("M134").Value = -62501304
$ws.Range("H136").Value = 9819.5
$ws.Range("I136").Value = 8729.25
$ws.Range("K136").Value = 26187.75
$ws.Range("M136").Value = -21087.75
$ws.Range("H137").Value = 4415.6523
$ws.Range("I137").Value = 4177.6
$ws.Range("J137").Value = 4862
$ws.Range("K137").Value = 12532.8
$ws.Range("L137").Value = 14586
$ws.Range("M137").Value = -7432.800000000001
$ws.Range("N137").Value = -24786
$ws.Range("H138").Value = 1827.125
$ws.Range("I138").Value = 1827.125
$ws.Range("K138").Value = 5481.375
$ws.Range("M138").Value = -341.375
$ws.Range("H139").Value = 3051.476
$ws.Range("I139").Value = 3328.6428
$ws.Range("K139").Value = 9985.928400000001
$ws.Range("M139").Value = -4845.928400000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2000.2778
$ws.Range("I3").Value = 1933.8667
$ws.Range("J3").Value = 2332.3333
$ws.Range("K3").Value = 1933.8667
$ws.Range("L3").Value = 2332.3333
$ws.Range("M3").Value = -1817.8667
$ws.Range("N3").Value = -2564.3333
$ws.Range("H14").Value = 700
$ws.Range("I14").Value = 700
$ws.Range("K14").Value = 700
$ws.Range("M14").Value = -532
$ws.Range("H80").Value = 3219.6
$ws.Range("I80").Value = 2824.5
$ws.Range("K80").Value = 2824.5
$ws.Range("M80").Value = -1826.5
$ws.Range("H83").Value = 3219.6
$ws.Range("I83").Value = 2824.5
$ws.Range("K83").Value = 14122.5
$ws.Range("M83").Value = -9130.5
$ws.Range("H97").Value = 1275.2632
$ws.Range("I97").Value = 1248.9412
$ws.Range("K97").Value = 1248.9412
$ws.Range("M97").Value = -752.9412
$ws.Range("H102").Value = 4396.4614
$ws.Range("I102").Value = 4019.25
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 4019.25
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -2397.25
$ws.Range("N102").Value = -8244

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3086.16
$ws.Range("I22").Value = 3265.923
$ws.Range("J22").Value = 2891.4167
$ws.Range("K22").Value = 3265.923
$ws.Range("L22").Value = 2891.4167
$ws.Range("M22").Value = -2970.923
$ws.Range("N22").Value = -3481.4167
$ws.Range("H27").Value = 3086.16
$ws.Range("I27").Value = 3265.923
$ws.Range("J27").Value = 2891.4167
$ws.Range("K27").Value = 3265.923
$ws.Range("L27").Value = 2891.4167
$ws.Range("M27").Value = -3158.923
$ws.Range("N27").Value = -3105.4167
$ws.Range("H46").Value = 1972.091
$ws.Range("J46").Value = 1750
$ws.Range("L46").Value = 1750
$ws.Range("N46").Value = -2126
$ws.Range("H127").Value = 99000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 17742.125
$ws.Range("J74").Value = 17742.125
$ws.Range("L74").Value = 17742.125
$ws.Range("N74").Value = -19614.125
$ws.Range("H77").Value = 17742.125
$ws.Range("J77").Value = 17742.125
$ws.Range("L77").Value = 53226.375
$ws.Range("N77").Value = -62586.375
$ws.Range("H107").Value = 31251012
$ws.Range("I107").Value = 35715376
$ws.Range("J107").Value = 478
$ws.Range("K107").Value = 107146128
$ws.Range("L107").Value = 1434
$ws.Range("M107").Value = -107144208
$ws.Range("N107").Value = -5274
